$wb = $excel.ActiveWorkbook

# New row (row 65) data for each of the 4 worksheets (MID_LFT_#1, MID_LFT_#2,
# MID_PLT_#1, MID_PLT_#2). Column layout:
#   A: time (datetime serial, formatted YYYY-MM-DD HH:MM:SS)
#   B: 总长            (hex string)
#   C: ID               (hex string)
#   D: 实际长度          (hex string)
#   E: 和校验            (hex string)
#   F: 总长_DEC          (number)
#   G: ID_DEC            (number)
#   H: 实际长度_DEC       (number)
#   I: 和校验_DEC         (number)

$rowsData = @(
    @{ Sheet = 1; B = "0x01,0x90"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"; D = "0x01,0x58"; E = "0x07"; F = 400; G = "5.68631262647113e+23"; H = 344; I = 7 },
    @{ Sheet = 2; B = "0x01,0x7c"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"; D = "0x01,0x54"; E = "0x19"; F = 380; G = "5.68432987514711e+23"; H = 340; I = 25 },
    @{ Sheet = 3; B = "0x00,0x6e"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"; D = "0x00,0x66"; E = "0x15"; F = 110; G = "5.68631262647113e+23"; H = 102; I = 15 },
    @{ Sheet = 4; B = "0x00,0x82"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"; D = "0x00,0x7B"; E = "0x9";  F = 130; G = "5.68631262647113e+23"; H = 123; I = 9 }
)

foreach ($rd in $rowsData) {
    $ws = $wb.Worksheets.Item($rd.Sheet)
    $r = 65

    $ws.Cells.Item($r, 1).Value = 45851.46114583333
    $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($r, 2).Value = $rd.B
    $ws.Cells.Item($r, 3).Value = $rd.C
    $ws.Cells.Item($r, 4).Value = $rd.D
    $ws.Cells.Item($r, 5).Value = $rd.E

    $ws.Cells.Item($r, 6).Value = $rd.F
    $ws.Cells.Item($r, 7).Value = [double]$rd.G
    $ws.Cells.Item($r, 8).Value = $rd.H
    $ws.Cells.Item($r, 9).Value = $rd.I
}
